$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily symbol-list refresh: updated Price (column D) and Volume(1h) (column E)
# figures for the affected coin rows. These cells hold literal text (e.g.
# "261.62", "0.42%"), not numeric/percentage values, so each write forces a
# text number format before assigning the value (otherwise Excel would infer
# a Number/Percentage type from strings that merely look numeric) and then
# restores the cell to the sheet's default "Normal" style so no stray
# formatting is introduced.
$updates = @{
    'D2'  = '261.62'
    'E2'  = '0.42%'
    'D3'  = '26.65'
    'E3'  = '-2.24%'
    'D4'  = '4.700'
    'E4'  = '0.33%'
    'E5'  = '-0.56%'
    'D6'  = '6.701'
    'E6'  = '0.67%'
    'D7'  = '0.8508'
    'E7'  = '-0.12%'
    'D8'  = '0.9114'
    'E8'  = '-1.14%'
    'D9'  = '0.1405'
    'E9'  = '0.37%'
    'D10' = '0.05110'
    'E10' = '7.30%'
    'D11' = '0.07088'
    'E11' = '0.02%'
    'D12' = '0.03113'
    'E12' = '2.01%'
    'D13' = '0.09032'
    'E13' = '-0.33%'
    'D14' = '0.001531'
    'E14' = '0.20%'
    'D15' = '0.0006188'
    'E15' = '1.74%'
    'D16' = '0.006004'
    'E16' = '-0.94%'
    'D17' = '3.449'
    'E17' = '-0.05%'
    'D18' = '3.168'
    'E18' = '0.66%'
    'E21' = '-0.76%'
    'D22' = '4.106'
    'E22' = '-0.11%'
    'D23' = '0.04227'
    'E23' = '0.07%'
    'D24' = '0.001178'
    'E24' = '-3.55%'
    'E25' = '6.93%'
    'E26' = '0.03%'
    'D27' = '0.0001938'
    'E27' = '23.04%'
    'D40' = '0.03952'
    'E40' = '2.53%'
    'E41' = '-0.13%'
    'D42' = '0.004173'
    'E42' = '1.99%'
    'D43' = '0.01391'
    'E43' = '-14.75%'
    'D44' = '0.002061'
    'E44' = '-7.03%'
    'D45' = '0.00005122'
    'E45' = '-0.75%'
    'D46' = '0.00000000750'
    'E46' = '0.03%'
    'D48' = '0.2579'
    'E48' = '61.92%'
    'D49' = '0.00002101'
    'E49' = '0.03%'
    'E50' = '0.03%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
